$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 16
$ws.Range("A8").Value = 37
$ws.Range("A10").Value = 39
$ws.Range("E10").Value = 473
$ws.Range("A11").Value = 59
$ws.Range("A12").Value = 5
$ws.Range("A14").Value = 28
$ws.Range("A15").Value = 20
$ws.Range("A16").Value = 17
$ws.Range("A20").Value = 41
$ws.Range("A21").Value = 31
$ws.Range("A22").Value = 19
$ws.Range("E22").Value = 528
$ws.Range("A25").Value = 58
$ws.Range("A27").Value = 8
$ws.Range("A28").Value = 36
$ws.Range("A33").Value = 34
$ws.Range("A36").Value = 6
$ws.Range("A38").Value = 35
$ws.Range("A39").Value = 50
$ws.Range("A41").Value = 27
$ws.Range("A42").Value = 40
$ws.Range("A43").Value = 32
$ws.Range("A44").Value = 38
$ws.Range("A48").Value = 18
$ws.Range("E48").Value = 845
$ws.Range("A49").Value = 49
$ws.Range("A52").Value = 4
$ws.Range("A53").Value = 9
$ws.Range("A56").Value = 30
$ws.Range("A58").Value = 7

$wb.Save()
